# Re-shuffle course placements across the weekday timetables ("mon",
# "tue", "wed", "thur", "fri" sheets) to match the new (if messier)
# arrangement described in the commit.  Each course occupies a
# contiguous 2-hour block (two adjacent cells in the same row), so the
# cell-by-cell writes below move / clear / create those blocks.

$wb = $excel.ActiveWorkbook

# --- Sheet: mon ---
$ws = $wb.Worksheets.Item("mon")
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("C4").Value = "CSC442"
$ws.Range("D4").Value = "CSC442"
$ws.Range("E20").Value = "MAT112"
$ws.Range("E4").Value = ""
$ws.Range("F20").Value = "MAT112"
$ws.Range("H21").Value = "CSC425"
$ws.Range("I13").Value = ""
$ws.Range("I21").Value = "CIT111"
$ws.Range("J13").Value = ""
$ws.Range("J21").Value = "CIT111"
$ws.Range("K18").Value = "CSC423"
$ws.Range("K21").Value = ""

# --- Sheet: tue ---
$ws = $wb.Worksheets.Item("tue")
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E15").Value = "CSC425"
$ws.Range("F15").Value = "CSC425"
$ws.Range("H17").Value = "CSC424"
$ws.Range("I20").Value = ""
$ws.Range("J15").Value = "BIO111"
$ws.Range("J17").Value = ""
$ws.Range("J20").Value = ""
$ws.Range("J21").Value = ""
$ws.Range("J25").Value = ""
$ws.Range("K17").Value = ""
$ws.Range("K20").Value = "CSC424"
$ws.Range("K25").Value = ""

# --- Sheet: wed ---
$ws = $wb.Worksheets.Item("wed")
$ws.Range("D20").Value = ""
$ws.Range("D21").Value = "TMC111"
$ws.Range("E20").Value = ""
$ws.Range("E21").Value = ""
$ws.Range("G2").Value = "CIS421"
$ws.Range("G25").Value = ""
$ws.Range("H13").Value = "CSC423"
$ws.Range("H2").Value = "CIS421"
$ws.Range("H20").Value = "CSC424"
$ws.Range("H25").Value = ""
$ws.Range("H8").Value = "BIO111"
$ws.Range("I11").Value = "CSC441"
$ws.Range("I13").Value = "CSC423"
$ws.Range("I20").Value = "CSC424"
$ws.Range("I4").Value = ""
$ws.Range("I8").Value = "BIO111"
$ws.Range("J11").Value = "CSC441"
$ws.Range("J20").Value = "DLD221"
$ws.Range("J4").Value = ""
$ws.Range("K20").Value = "DLD221"

# --- Sheet: thur ---
$ws = $wb.Worksheets.Item("thur")
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("F21").Value = "MAT111"
$ws.Range("F4").Value = "CSC111"
$ws.Range("G17").Value = "CSC424"
$ws.Range("G21").Value = "MAT111"
$ws.Range("H17").Value = "CSC424"
$ws.Range("I20").Value = ""
$ws.Range("J11").Value = ""
$ws.Range("J21").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("K18").Value = "CSC442"
$ws.Range("K20").Value = "CHM111"

# --- Sheet: fri ---
$ws = $wb.Worksheets.Item("fri")
$ws.Range("C21").Value = ""
$ws.Range("E20").Value = ""
